$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 191.4762
$ws.Range("I33").Value = 148.47368
$ws.Range("K33").Value = 148.47368
$ws.Range("M33").Value = 80.52632
$ws.Range("H98").Value = 1947.4231
$ws.Range("I98").Value = 1909.4286
$ws.Range("J98").Value = 2107
$ws.Range("K98").Value = 1909.4286
$ws.Range("L98").Value = 2107
$ws.Range("M98").Value = -411.4286
$ws.Range("N98").Value = -5103
$ws.Range("H100").Value = 3070.16
$ws.Range("I100").Value = 1525
$ws.Range("J100").Value = 5036.727
$ws.Range("K100").Value = 1525
$ws.Range("L100").Value = 5036.727
$ws.Range("M100").Value = -984
$ws.Range("N100").Value = -6118.727
$ws.Range("H122").Value = 1947.4231
$ws.Range("I122").Value = 1909.4286
$ws.Range("J122").Value = 2107
$ws.Range("K122").Value = 5728.2858
$ws.Range("L122").Value = 6321
$ws.Range("M122").Value = -3278.2858
$ws.Range("N122").Value = -11221
$ws.Range("H132").Value = 7629.069
$ws.Range("I132").Value = 7629.069
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 22887.207
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -20357.207
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 3509.132
$ws.Range("I138").Value = 2804.95
$ws.Range("K138").Value = 8414.849999999999
$ws.Range("M138").Value = -3274.849999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1266090.4
$ws.Range("J32").Value = 51249.5
$ws.Range("L32").Value = 51249.5
$ws.Range("N32").Value = -51823.5
$ws.Range("H97").Value = 673.1667
$ws.Range("I97").Value = 673.1667
$ws.Range("K97").Value = 673.1667
$ws.Range("M97").Value = -177.1667
$ws.Range("H122").Value = 2208
$ws.Range("I122").Value = 2052.4
$ws.Range("J122").Value = 2597
$ws.Range("K122").Value = 6157.200000000001
$ws.Range("L122").Value = 7791
$ws.Range("M122").Value = -3707.200000000001
$ws.Range("N122").Value = -12691
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 35629.668
$ws.Range("I20").Value = 40347.027
$ws.Range("K20").Value = 40347.027
$ws.Range("M20").Value = -40100.027
$ws.Range("H86").Value = 6132
$ws.Range("I86").Value = 3548
$ws.Range("J86").Value = 6778
$ws.Range("K86").Value = 3548
$ws.Range("L86").Value = 6778
$ws.Range("M86").Value = -2425
$ws.Range("N86").Value = -9024
$ws.Range("H89").Value = 6132
$ws.Range("I89").Value = 3548
$ws.Range("J89").Value = 6778
$ws.Range("K89").Value = 17740
$ws.Range("L89").Value = 33890
$ws.Range("M89").Value = -12124
$ws.Range("N89").Value = -45122
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 1514.25
$ws.Range("I107").Value = 1135.3572
$ws.Range("K107").Value = 1135.3572
$ws.Range("M107").Value = 784.6428000000001
$ws.Range("H134").Value = 2977579.2
$ws.Range("I134").Value = 1457.9231
$ws.Range("K134").Value = 4373.7693
$ws.Range("M134").Value = -1838.7693

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2478.875
$ws.Range("I99").Value = 2089.4
$ws.Range("J99").Value = 3128
$ws.Range("K99").Value = 2089.4
$ws.Range("L99").Value = 3128
$ws.Range("M99").Value = -591.4000000000001
$ws.Range("N99").Value = -6124
$ws.Range("H126").Value = 2478.875
$ws.Range("I126").Value = 2089.4
$ws.Range("J126").Value = 3128
$ws.Range("K126").Value = 6268.200000000001
$ws.Range("L126").Value = 9384
$ws.Range("M126").Value = -3798.200000000001
$ws.Range("N126").Value = -14324
$ws.Range("H132").Value = 2968.35
$ws.Range("I132").Value = 2775.1667
$ws.Range("K132").Value = 8325.500100000001
$ws.Range("M132").Value = -5795.500100000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 303.35715
$ws.Range("J23").Value = 294.5
$ws.Range("L23").Value = 883.5
$ws.Range("N23").Value = -1353.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 39666.668
$ws.Range("I94").Value = 30000
$ws.Range("K94").Value = 30000
$ws.Range("M94").Value = -29324
$ws.Range("H132").Value = 25107
$ws.Range("I132").Value = 13045.363
$ws.Range("J132").Value = 69333
$ws.Range("K132").Value = 39136.089
$ws.Range("L132").Value = 207999
$ws.Range("M132").Value = -36606.089
$ws.Range("N132").Value = -213059

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3149.6191
$ws.Range("I22").Value = 2856.2856
$ws.Range("J22").Value = 3296.2856
$ws.Range("K22").Value = 2856.2856
$ws.Range("L22").Value = 3296.2856
$ws.Range("M22").Value = -2561.2856
$ws.Range("N22").Value = -3886.2856
$ws.Range("H27").Value = 3149.6191
$ws.Range("I27").Value = 2856.2856
$ws.Range("J27").Value = 3296.2856
$ws.Range("K27").Value = 2856.2856
$ws.Range("L27").Value = 3296.2856
$ws.Range("M27").Value = -2749.2856
$ws.Range("N27").Value = -3510.2856
$ws.Range("H40").Value = 6195.3
$ws.Range("I40").Value = 5696.5
$ws.Range("J40").Value = 6943.5
$ws.Range("K40").Value = 5696.5
$ws.Range("L40").Value = 6943.5
$ws.Range("M40").Value = -5560.5
$ws.Range("N40").Value = -7215.5
$ws.Range("H55").Value = 1552.258
$ws.Range("J55").Value = 1428.9
$ws.Range("L55").Value = 1428.9
$ws.Range("N55").Value = -1774.9
$ws.Range("H68").Value = 2956.08
$ws.Range("I68").Value = 2333.5557
$ws.Range("J68").Value = 4556.857
$ws.Range("K68").Value = 2333.5557
$ws.Range("L68").Value = 4556.857
$ws.Range("M68").Value = -1584.5557
$ws.Range("N68").Value = -6054.857
$ws.Range("H71").Value = 2956.08
$ws.Range("I71").Value = 2333.5557
$ws.Range("J71").Value = 4556.857
$ws.Range("K71").Value = 11667.7785
$ws.Range("L71").Value = 22784.285
$ws.Range("M71").Value = -7923.7785
$ws.Range("N71").Value = -30272.285
$ws.Range("H93").Value = 3029.261
$ws.Range("I93").Value = 2173.8948
$ws.Range("J93").Value = 7092.25
$ws.Range("K93").Value = 2173.8948
$ws.Range("L93").Value = 7092.25
$ws.Range("M93").Value = -925.8948
$ws.Range("N93").Value = -9588.25
$ws.Range("H122").Value = 3302.15
$ws.Range("I122").Value = 2938.6667
$ws.Range("K122").Value = 8816.000100000001
$ws.Range("M122").Value = -6366.000100000001
$ws.Range("H132").Value = 1518659.6
$ws.Range("I132").Value = 3336811.8
$ws.Range("J132").Value = 3532.9167
$ws.Range("K132").Value = 10010435.4
$ws.Range("L132").Value = 10598.7501
$ws.Range("M132").Value = -10007905.4
$ws.Range("N132").Value = -15658.7501
$ws.Range("H136").Value = 16669039
$ws.Range("I136").Value = 10419092
$ws.Range("J136").Value = 41668824
$ws.Range("K136").Value = 31257276
$ws.Range("L136").Value = 125006472
$ws.Range("M136").Value = -31254726
$ws.Range("N136").Value = -125011572

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 33263.332
$ws.Range("I62").Value = 33263.332
$ws.Range("K62").Value = 33263.332
$ws.Range("M62").Value = -32639.332
$ws.Range("H65").Value = 33263.332
$ws.Range("I65").Value = 33263.332
$ws.Range("K65").Value = 166316.66
$ws.Range("M65").Value = -163196.66
$ws.Range("H113").Value = 1553.4412
$ws.Range("I113").Value = 1314
$ws.Range("K113").Value = 3942
$ws.Range("M113").Value = -1772
$ws.Range("H132").Value = 12823357
$ws.Range("I132").Value = 18521160
$ws.Range("J132").Value = 3299.75
$ws.Range("K132").Value = 55563480
$ws.Range("L132").Value = 9899.25
$ws.Range("M132").Value = -55560950
$ws.Range("N132").Value = -14959.25

Write-Host "Done applying Brynhildr_Profits updates"